$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Birth/death years for rows 122-166 (Expresionismo Abstracto .. Arte Conceptual / Robert Morris) ---
    $ws.Range("C122").Value = 1899
    $ws.Range("D122").Value = 1968
    $ws.Range("C123").Value = 1963
    $ws.Range("C124").Value = 1928
    $ws.Range("D124").Value = 2011
    $ws.Range("C125").Value = 1888
    $ws.Range("D125").Value = 1976
    $ws.Range("C126").Value = 1930
    $ws.Range("C127").Value = 1925
    $ws.Range("D127").Value = 2008
    $ws.Range("C128").Value = 1929
    $ws.Range("C129").Value = 1939
    $ws.Range("C130").Value = 1920
    $ws.Range("C131").Value = 1928
    $ws.Range("D131").Value = 1987
    $ws.Range("C132").Value = 1931
    $ws.Range("D132").Value = 2004
    $ws.Range("C133").Value = 1937
    $ws.Range("C134").Value = 1933
    $ws.Range("D134").Value = 2017
    $ws.Range("C135").Value = 1901
    $ws.Range("D135").Value = 1978
    $ws.Range("C136").Value = 1923
    $ws.Range("D136").Value = 1997
    $ws.Range("C137").Value = 1941
    $ws.Range("D137").Value = 2010
    $ws.Range("C138").Value = 1924
    $ws.Range("D138").Value = 2000
    $ws.Range("C139").Value = 1937
    $ws.Range("C140").Value = 1922
    $ws.Range("D140").Value = 2011
    $ws.Range("C141").Value = 1938
    $ws.Range("C142").Value = 1932
    $ws.Range("C143").Value = 1909
    $ws.Range("D143").Value = 1992
    $ws.Range("C144").Value = 1932
    $ws.Range("C145").Value = 1935
    $ws.Range("C146").Value = 1940
    $ws.Range("C147").Value = 1911
    $ws.Range("D147").Value = 2010
    $ws.Range("C148").Value = 1930
    $ws.Range("D148").Value = 1998
    $ws.Range("C149").Value = 1933
    $ws.Range("D149").Value = 1996
    $ws.Range("C150").Value = 1928
    $ws.Range("D150").Value = 2007
    $ws.Range("C151").Value = 1930
    $ws.Range("C152").Value = 1936
    $ws.Range("D152").Value = 1970
    $ws.Range("C153").Value = 1935
    $ws.Range("C154").Value = 1928
    $ws.Range("D154").Value = 1994
    $ws.Range("C155").Value = 1939
    $ws.Range("C156").Value = 1921
    $ws.Range("D156").Value = 1986
    $ws.Range("C157").Value = 1938
    $ws.Range("C158").Value = 1945
    $ws.Range("C159").Value = 1941
    $ws.Range("C160").Value = 1938
    $ws.Range("D160").Value = 1973
    $ws.Range("C161").Value = 1940
    $ws.Range("D161").Value = 1994
    $ws.Range("C162").Value = 1924
    $ws.Range("D162").Value = 1976
    $ws.Range("C163").Value = 1933
    $ws.Range("D163").Value = 2014
    $ws.Range("C164").Value = 1932
    $ws.Range("D164").Value = 2006
    $ws.Range("C165").Value = 1943
    $ws.Range("D165").Value = 1978
    $ws.Range("C166").Value = 1931
    $ws.Range("D166").Value = 2018

# --- Split "Christo and Jeanne-Claude" (row 167) into two separate artist rows ---
$ws.Rows.Item(168).Insert()

$ws.Range("B167").Value = "Christo Vladimirov Javacheff"
$ws.Range("C167").Value = 1935
$ws.Range("D167").Value = 2020

$ws.Range("A168").Value = "Arte Conceptual"
$ws.Range("B168").Value = "Jeanne-Claude Denat de Guillebon"
$ws.Range("C168").Value = 1935
$ws.Range("D168").Value = 2009

# --- Birth/death years for the remaining rows (now shifted down by one) ---
    $ws.Range("C169").Value = 1943
    $ws.Range("C170").Value = 1935
    $ws.Range("D170").Value = 2013
    $ws.Range("C171").Value = 1954
    $ws.Range("C172").Value = 1949
    $ws.Range("C173").Value = 1951
    $ws.Range("C174").Value = 1960
    $ws.Range("D174").Value = 1988
    $ws.Range("C175").Value = 1948
    $ws.Range("C176").Value = 1945
    $ws.Range("C177").Value = 1958
    $ws.Range("D177").Value = 1990
    $ws.Range("C178").Value = 1945
    $ws.Range("C179").Value = 1954
    $ws.Range("D179").Value = 2012
    $ws.Range("C180").Value = 1952
    $ws.Range("D180").Value = 2013
    $ws.Range("C181").Value = 1953
    $ws.Range("D181").Value = 1997
    $ws.Range("C182").Value = 1955
    $ws.Range("C183").Value = 1965
    $ws.Range("C184").Value = 1954
    $ws.Range("C185").Value = 1945
    $ws.Range("C186").Value = 1954
    $ws.Range("C187").Value = 1955
    $ws.Range("C188").Value = 1957
    $ws.Range("C189").Value = 1960
    $ws.Range("C190").Value = 1954
    $ws.Range("C191").Value = 1959
    $ws.Range("C192").Value = 1960
    $ws.Range("C193").Value = 1955
    $ws.Range("C194").Value = 1966

# --- Cosmetic refresh: column width, selection and view position ---
$ws.Columns.Item(2).ColumnWidth = 32.28515625
$ws.Activate()
$ws.Range("C200").Select()
